# Set the "Status" value to "Sent" for the rows that don't already have it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 5, 6, 7, 10)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "Sent"
}
